# Update the "F" column (想去人数 / "want to go" count) values across all four
# worksheets to reflect the regenerated site data (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 470  # was 469
$ws.Range("F5").Value = 8831  # was 8824
$ws.Range("F6").Value = 18  # was 17
$ws.Range("F7").Value = 11318  # was 11312
$ws.Range("F20").Value = 425  # was 424
$ws.Range("F22").Value = 729  # was 728
$ws.Range("F23").Value = 656  # was 654
$ws.Range("F24").Value = 365  # was 364
$ws.Range("F29").Value = 1350  # was 1343
$ws.Range("F30").Value = 31  # was 30
$ws.Range("F34").Value = 1427  # was 1426
$ws.Range("F36").Value = 312  # was 237
$ws.Range("F37").Value = 23  # was 15
$ws.Range("F38").Value = 359  # was 358
$ws.Range("F39").Value = 333  # was 330
$ws.Range("F41").Value = 149  # was 148
$ws.Range("F43").Value = 399  # was 398
$ws.Range("F45").Value = 821  # was 819
$ws.Range("F48").Value = 176  # was 172
$ws.Range("F49").Value = 161  # was 160

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F18").Value = 71  # was 70
$ws.Range("F19").Value = 113  # was 112
$ws.Range("F24").Value = 80  # was 79
$ws.Range("F25").Value = 399  # was 398

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 2843  # was 2841
$ws.Range("F4").Value = 353  # was 351
$ws.Range("F5").Value = 218  # was 217

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 353  # was 351
$ws.Range("F6").Value = 218  # was 217
$ws.Range("F7").Value = 8831  # was 8824
$ws.Range("F8").Value = 18  # was 17
$ws.Range("F9").Value = 11318  # was 11313
$ws.Range("F18").Value = 425  # was 424
$ws.Range("F20").Value = 729  # was 728
$ws.Range("F21").Value = 656  # was 654
$ws.Range("F22").Value = 365  # was 364
$ws.Range("F29").Value = 1350  # was 1343
$ws.Range("F30").Value = 31  # was 30
$ws.Range("F35").Value = 1427  # was 1426
$ws.Range("F36").Value = 312  # was 238
$ws.Range("F37").Value = 359  # was 358
$ws.Range("F39").Value = 149  # was 148
$ws.Range("F41").Value = 399  # was 398
$ws.Range("F45").Value = 399  # was 398
$ws.Range("F48").Value = 176  # was 172
$ws.Range("F49").Value = 161  # was 160

